$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1877.4445  # was 1899.7
$ws.Range("J40").Value = 2099.3333  # was 2099.5
$ws.Range("L40").Value = 2099.3333  # was 2099.5
$ws.Range("N40").Value = -2449.3333  # was -2449.5
$ws.Range("H53").Value = 607  # was 567.0909
$ws.Range("I53").Value = 525  # was 511.125
$ws.Range("J53").Value = 894  # was 716.3333
$ws.Range("K53").Value = 525  # was 511.125
$ws.Range("L53").Value = 894  # was 716.3333
$ws.Range("M53").Value = 112  # was 125.875
$ws.Range("N53").Value = -2168  # was -1990.3333
$ws.Range("H64").Value = 2600  # was 0
$ws.Range("I64").Value = 2600  # was 0
$ws.Range("K64").Value = 2600  # was 0
$ws.Range("M64").Value = -2352
$ws.Range("H67").Value = 2600  # was 0
$ws.Range("I67").Value = 2600  # was 0
$ws.Range("K67").Value = 2600  # was 0
$ws.Range("M67").Value = -1742
$ws.Range("I137").Value = 1000  # was 0
$ws.Range("K137").Value = 3000  # was 0
$ws.Range("M137").Value = -450

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4515.1113  # was 2190.1538
$ws.Range("I32").Value = 1684.8  # was 2190.1538
$ws.Range("J32").Value = 18666.666  # was 0
$ws.Range("K32").Value = 1684.8  # was 2190.1538
$ws.Range("L32").Value = 18666.666  # was 0
$ws.Range("M32").Value = -1397.8  # was -1903.1538
$ws.Range("N32").Value = -19240.666
$ws.Range("H45").Value = 2884.75  # was 2916.8
$ws.Range("I45").Value = 2884.75  # was 2916.8
$ws.Range("K45").Value = 2884.75  # was 2916.8
$ws.Range("M45").Value = -2507.75  # was -2539.8
$ws.Range("H80").Value = 83748  # was 100108
$ws.Range("I80").Value = 87501  # was 0
$ws.Range("J80").Value = 79995  # was 100108
$ws.Range("K80").Value = 87501  # was 0
$ws.Range("L80").Value = 79995  # was 100108
$ws.Range("M80").Value = -86503
$ws.Range("N80").Value = -81991  # was -102104
$ws.Range("H83").Value = 83748  # was 100108
$ws.Range("I83").Value = 87501  # was 0
$ws.Range("J83").Value = 79995  # was 100108
$ws.Range("K83").Value = 262503  # was 0
$ws.Range("L83").Value = 239985  # was 300324
$ws.Range("M83").Value = -257511
$ws.Range("N83").Value = -249969  # was -310308
$ws.Range("H97").Value = 1537.1428  # was 1628.5385
$ws.Range("I97").Value = 1022.8  # was 1097.6666
$ws.Range("K97").Value = 1022.8  # was 1097.6666
$ws.Range("M97").Value = -526.8  # was -601.6666
$ws.Range("H122").Value = 4579.75  # was 5766.6665
$ws.Range("I122").Value = 4579.75  # was 5766.6665
$ws.Range("K122").Value = 13739.25  # was 17299.9995
$ws.Range("M122").Value = -11289.25  # was -14849.9995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1551.3334  # was 1037.7778
$ws.Range("I20").Value = 1374.75  # was 1037.7778
$ws.Range("J20").Value = 1904.5  # was 0
$ws.Range("K20").Value = 1374.75  # was 1037.7778
$ws.Range("L20").Value = 1904.5  # was 0
$ws.Range("M20").Value = -1127.75  # was -790.7778000000001
$ws.Range("N20").Value = -2398.5
$ws.Range("H105").Value = 1483.3334  # was 5000
$ws.Range("I105").Value = 1483.3334  # was 5000
$ws.Range("K105").Value = 1483.3334  # was 5000
$ws.Range("M105").Value = 263.6666  # was -3253

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 5507.4443  # was 5495.778
$ws.Range("I11").Value = 1005  # was 552.5
$ws.Range("J11").Value = 6070.25  # was 6908.143
$ws.Range("K11").Value = 1005  # was 552.5
$ws.Range("L11").Value = 6070.25  # was 6908.143
$ws.Range("M11").Value = -865  # was -412.5
$ws.Range("N11").Value = -6350.25  # was -7188.143
$ws.Range("H22").Value = 1133.3334  # was 1499.5
$ws.Range("I22").Value = 700  # was 999
$ws.Range("K22").Value = 700  # was 999
$ws.Range("M22").Value = -350  # was -649
$ws.Range("H31").Value = 2200  # was 1945.5
$ws.Range("I31").Value = 2200  # was 1945.5
$ws.Range("K31").Value = 2200  # was 1945.5
$ws.Range("M31").Value = -1905  # was -1650.5
$ws.Range("H34").Value = 2200  # was 1945.5
$ws.Range("I34").Value = 2200  # was 1945.5
$ws.Range("K34").Value = 2200  # was 1945.5
$ws.Range("M34").Value = -1998  # was -1743.5
$ws.Range("H43").Value = 71249.75  # was 71250
$ws.Range("J43").Value = 71249.75  # was 71250
$ws.Range("L43").Value = 71249.75  # was 71250
$ws.Range("N43").Value = -71617.75  # was -71618
$ws.Range("H101").Value = 71249.75  # was 71250
$ws.Range("J101").Value = 71249.75  # was 71250
$ws.Range("L101").Value = 71249.75  # was 71250
$ws.Range("N101").Value = -77739.75  # was -77740
$ws.Range("H103").Value = 50000  # was 44975
$ws.Range("I103").Value = 50000  # was 44975
$ws.Range("K103").Value = 50000  # was 44975
$ws.Range("M103").Value = -48828  # was -43803
$ws.Range("H122").Value = 18918  # was 12899.5
$ws.Range("J122").Value = 24957.334  # was 24999
$ws.Range("L122").Value = 74872.00199999999  # was 74997
$ws.Range("N122").Value = -79772.00199999999  # was -79897
$ws.Range("H134").Value = 1010.1667  # was 995.125
$ws.Range("I134").Value = 852.2  # was 826.8333
$ws.Range("J134").Value = 1800  # was 1500
$ws.Range("K134").Value = 2556.6  # was 2480.4999
$ws.Range("L134").Value = 5400  # was 4500
$ws.Range("M134").Value = -21.60000000000036  # was 54.5001000000002
$ws.Range("N134").Value = -10470  # was -9570
$ws.Range("H140").Value = 0  # was 95000
$ws.Range("J140").Value = 0  # was 95000
$ws.Range("L140").Value = 0  # was 95000
$ws.Range("N140").ClearContents()  # was -105360

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 96.333336  # was 650.75
$ws.Range("I11").Value = 105.8  # was 220
$ws.Range("J11").Value = 49  # was 1512.25
$ws.Range("K11").Value = 317.4  # was 660
$ws.Range("L11").Value = 147  # was 4536.75
$ws.Range("M11").Value = -177.4  # was -520
$ws.Range("N11").Value = -427  # was -4816.75
$ws.Range("H26").Value = 1386  # was 1397.25
$ws.Range("I26").Value = 97.5  # was 547.5
$ws.Range("J26").Value = 2245  # was 2247
$ws.Range("K26").Value = 292.5  # was 1642.5
$ws.Range("L26").Value = 6735  # was 6741
$ws.Range("M26").Value = -4.5  # was -1354.5
$ws.Range("N26").Value = -7311  # was -7317
$ws.Range("H93").Value = 2000  # was 0
$ws.Range("J93").Value = 2000  # was 0
$ws.Range("L93").Value = 6000  # was 0
$ws.Range("N93").Value = -9744
$ws.Range("H105").Value = 10000  # was 9514.5
$ws.Range("J105").Value = 10000  # was 9514.5
$ws.Range("L105").Value = 30000  # was 28543.5
$ws.Range("N105").Value = -35242  # was -33785.5
$ws.Range("H109").Value = 659  # was 469
$ws.Range("I109").Value = 690.8  # was 469
$ws.Range("J109").Value = 500  # was 0
$ws.Range("K109").Value = 2072.4  # was 1407
$ws.Range("L109").Value = 1500  # was 0
$ws.Range("M109").Value = -1032.4  # was -367
$ws.Range("N109").Value = -3580
$ws.Range("H121").Value = 395.7143  # was 365
$ws.Range("I121").Value = 295  # was 365
$ws.Range("J121").Value = 1000  # was 0
$ws.Range("K121").Value = 885  # was 1095
$ws.Range("L121").Value = 3000  # was 0
$ws.Range("M121").Value = 425  # was 215
$ws.Range("N121").Value = -5620
$ws.Range("H124").Value = 2250  # was 2500
$ws.Range("I124").Value = 2250  # was 2500
$ws.Range("K124").Value = 6750  # was 7500
$ws.Range("M124").Value = -1840  # was -2590
$ws.Range("H132").Value = 824.6667  # was 997
$ws.Range("J132").Value = 480  # was 0
$ws.Range("L132").Value = 4320  # was 0
$ws.Range("N132").Value = -9380
$ws.Range("H140").Value = 1457.5  # was 1250
$ws.Range("I140").Value = 1457.5  # was 1250
$ws.Range("K140").Value = 4372.5  # was 3750
$ws.Range("M140").Value = 807.5  # was 1430

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0  # was 6499
$ws.Range("J70").Value = 0  # was 6499
$ws.Range("L70").Value = 0  # was 6499
$ws.Range("N70").ClearContents()  # was -7039
$ws.Range("H73").Value = 0  # was 6499
$ws.Range("J73").Value = 0  # was 6499
$ws.Range("L73").Value = 0  # was 6499
$ws.Range("N73").ClearContents()  # was -8371
$ws.Range("H122").Value = 3888.8235  # was 4079.1875
$ws.Range("I122").Value = 2887.0833  # was 3072.9092
$ws.Range("K122").Value = 8661.249899999999  # was 9218.7276
$ws.Range("M122").Value = -6211.249899999999  # was -6768.7276

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5000  # was 0
$ws.Range("I7").Value = 5000  # was 0
$ws.Range("K7").Value = 5000  # was 0
$ws.Range("M7").Value = -4888
$ws.Range("H40").Value = 1015999.8  # was 1016000
$ws.Range("J40").Value = 1676666.4  # was 1676666.6
$ws.Range("L40").Value = 1676666.4  # was 1676666.6
$ws.Range("N40").Value = -1676938.4  # was -1676938.6
$ws.Range("H126").Value = 5000  # was 0
$ws.Range("I126").Value = 5000  # was 0
$ws.Range("K126").Value = 15000  # was 0
$ws.Range("M126").Value = -12530
$ws.Range("H136").Value = 834208.3  # was 1000870
$ws.Range("I136").Value = 834208.3  # was 1000870
$ws.Range("K136").Value = 2502624.9  # was 3002610
$ws.Range("M136").Value = -2500074.9  # was -3000060

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 10400  # was 10333.333
$ws.Range("H132").Value = 1295.9166  # was 1346.0834
$ws.Range("I132").Value = 1105.1  # was 1195.7273
$ws.Range("J132").Value = 2250  # was 3000
$ws.Range("K132").Value = 3315.3  # was 3587.1819
$ws.Range("L132").Value = 6750  # was 9000
$ws.Range("M132").Value = -785.2999999999997  # was -1057.1819
$ws.Range("N132").Value = -11810  # was -14060
$ws.Range("H136").Value = 1959.4  # was 1969.3
$ws.Range("I136").Value = 1959.4  # was 1969.3
$ws.Range("K136").Value = 5878.200000000001  # was 5907.9
$ws.Range("M136").Value = -3328.200000000001  # was -3357.9
